$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
  "Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404",
  "Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404",
  "diff",
  "Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410",
  "Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U85"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = $null

Write-Host "done"
